$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free approach: D-column numeric-looking text values are written with a
# leading apostrophe so Excel keeps them as text (matching source data that uses
# dotted/grouped number formats like "29.163.29"), then the style is reset to
# Normal so no stray quote-prefix formatting is left behind.

# --- Row swaps: Polkadot/Polygon (rows 15-16) and Stellar/EthereumClassic (rows 26-27) ---
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.523"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.70%  "

$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "'3.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.35%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'15.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.69%  "

$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "'0.107"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.99%  "

# --- Price / Volume updates ---
$ws.Range("D2").Value = "'29.163.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.00%  "
$ws.Range("D3").Value = "'1.578.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.76%  "
$ws.Range("D4").Value = "'0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").Value = "'211.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").Value = "'0.512"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.94%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "'26.31"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +11.05%  "
$ws.Range("E9").Value = "  +2.19%  "
$ws.Range("D10").Value = "'0.0592"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.56%  "
$ws.Range("E11").Value = "  +1.49%  "
$ws.Range("D12").Value = "'1.804.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.74%  "
$ws.Range("D13").Value = "'1.577.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.69%  "
$ws.Range("D14").Value = "'29.195.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.16%  "
$ws.Range("D17").Value = "'62.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.93%  "
$ws.Range("D18").Value = "'237.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.20%  "
$ws.Range("D19").Value = "'7.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0690"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.46%  "
$ws.Range("D21").Value = "'0.994"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("D22").Value = "'3.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.88%  "
$ws.Range("D23").Value = "'9.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.50%  "
$ws.Range("E24").Value = "  +4.74%  "
$ws.Range("D25").Value = "'153.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.56%  "
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").Value = "'0.0466"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("E32").Value = "  +1.52%  "
$ws.Range("D33").Value = "'1.424.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.67%  "
$ws.Range("D34").Value = "'3.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("D35").Value = "'1.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.73%  "
$ws.Range("E36").Value = "  +0.99%  "
$ws.Range("D37").Value = "'2.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.86%  "
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("E39").Value = "  +0.99%  "
$ws.Range("E40").Value = "  +3.01%  "
$ws.Range("D41").Value = "'1.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.95%  "
$ws.Range("D42").Value = "'52.87"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +26.53%  "
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("D44").Value = "'0.789"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.44%  "
$ws.Range("D45").Value = "'0.0469"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("D46").Value = "'64.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.16%  "
$ws.Range("D47").Value = "'5.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("D48").Value = "'1.716.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.76%  "
$ws.Range("E49").Value = "  -6.47%  "
$ws.Range("D50").Value = "'85.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").Value = "'0.0₆0103"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.39%  "
